$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number (45174 = 2023-09-05) for
# every data row (2..176). Bump it forward one day to 45175 (2023-09-06).
$ws.Range("C2:C176").Value = 45175
